# Apply crypto list update (prices, 1h volume %, and two row swaps)
# Commit: Updated cryptos list on Thu Dec 14 17:17:56 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $r = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "6.00", "0.650")
    # keep their exact original formatting instead of being parsed as numbers.
    $r.NumberFormat = "@"
    $r.Value = $Text
    # Restore default cell style so no stray formatting is introduced.
    $r.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "42.781.88"
Set-TextValue "E2" "  +2.08%  "
# Row 3
Set-TextValue "D3" "2.294.59"
Set-TextValue "E3" "  +3.91%  "
# Row 4
Set-TextValue "E4" "  -0.01%  "
# Row 5
Set-TextValue "D5" "251.90"
Set-TextValue "E5" "  -0.04%  "
# Row 6
Set-TextValue "D6" "0.639"
Set-TextValue "E6" "  +4.17%  "
# Row 7
Set-TextValue "D7" "74.51"
Set-TextValue "E7" "  +9.98%  "
# Row 8
Set-TextValue "E8" "  -0.07%  "
# Row 9
Set-TextValue "D9" "0.650"
Set-TextValue "E9" "  +4.88%  "
# Row 10
Set-TextValue "D10" "39.76"
Set-TextValue "E10" "  +2.04%  "
# Row 11
Set-TextValue "D11" "0.0988"
Set-TextValue "E11" "  +5.42%  "
# Row 12
Set-TextValue "D12" "59.16"
Set-TextValue "E12" "  -0.47%  "
# Row 13
Set-TextValue "D13" "7.38"
Set-TextValue "E13" "  +4.84%  "
# Row 14
Set-TextValue "E14" "  +1.50%  "
# Row 15
Set-TextValue "D15" "2.632.09"
Set-TextValue "E15" "  +3.71%  "
# Row 16
Set-TextValue "D16" "15.37"
Set-TextValue "E16" "  +6.51%  "
# Row 17
Set-TextValue "D17" "0.882"
Set-TextValue "E17" "  +1.34%  "
# Row 18
Set-TextValue "D18" "2.282.86"
Set-TextValue "E18" "  +1.51%  "
# Row 19
Set-TextValue "D19" "42.705.27"
Set-TextValue "E19" "  +2.06%  "
# Row 20
Set-TextValue "E20" "  +4.80%  "
# Row 21
Set-TextValue "D21" "6.32"
Set-TextValue "E21" "  +3.17%  "
# Row 22
Set-TextValue "D22" "72.51"
Set-TextValue "E22" "  +0.27%  "
# Row 23
Set-TextValue "D23" "233.31"
Set-TextValue "E23" "  +0.87%  "
# Row 24
Set-TextValue "E24" "  +10.30%  "
# Row 25
Set-TextValue "D25" "3.91"
Set-TextValue "E25" "  +0.63%  "
# Row 26
Set-TextValue "D26" "11.61"
Set-TextValue "E26" "  +4.30%  "
# Row 27
Set-TextValue "E27" "  -0.26%  "
# Row 28
Set-TextValue "D28" "2.42"
Set-TextValue "E28" "  +0.33%  "
# Row 29
Set-TextValue "D29" "3.63"
Set-TextValue "E29" "  -1.48%  "
# Row 30
Set-TextValue "D30" "2.19"
Set-TextValue "E30" "  +1.27%  "
# Row 31
Set-TextValue "D31" "167.14"
Set-TextValue "E31" "  +0.43%  "
# Row 32
Set-TextValue "D32" "21.15"
Set-TextValue "E32" "  +3.70%  "
# Row 33
Set-TextValue "D33" "6.43"
Set-TextValue "E33" "  +9.77%  "
# Row 34
Set-TextValue "E34" "  +5.29%  "
# Row 35
Set-TextValue "B35" "Hedera"
Set-TextValue "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.0823"
Set-TextValue "E35" "  +5.47%  "
# Row 36
Set-TextValue "B36" "InjectiveProtocol"
Set-TextValue "C36" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D36" "32.32"
Set-TextValue "E36" "  +25.88%  "
# Row 37
Set-TextValue "D37" "0.127"
Set-TextValue "E37" "  +4.21%  "
# Row 38
Set-TextValue "D38" "4.75"
Set-TextValue "E38" "  +17.03%  "
# Row 39
Set-TextValue "D39" "4.78"
Set-TextValue "E39" "  +4.18%  "
# Row 40
Set-TextValue "D40" "0.0308"
Set-TextValue "E40" "  +0.18%  "
# Row 41
Set-TextValue "D41" "14.28"
Set-TextValue "E41" "  +19.00%  "
# Row 42
Set-TextValue "D42" "2.35"
Set-TextValue "E42" "  +5.86%  "
# Row 43
Set-TextValue "D43" "6.00"
Set-TextValue "E43" "  +6.58%  "
# Row 44
Set-TextValue "D44" "0.214"
Set-TextValue "E44" "  +9.50%  "
# Row 45
Set-TextValue "B45" "MultiversX"
Set-TextValue "C45" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D45" "62.29"
Set-TextValue "E45" "  +1.65%  "
# Row 46
Set-TextValue "B46" "FraxShare"
Set-TextValue "C46" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "9.15"
Set-TextValue "E46" "  +6.76%  "
# Row 47
Set-TextValue "D47" "4.87"
# Row 48
Set-TextValue "D48" "0.103"
Set-TextValue "E48" "  +3.89%  "
# Row 49
Set-TextValue "E49" "  +0.43%  "
# Row 50
Set-TextValue "E50" "  +3.27%  "
# Row 51
Set-TextValue "D51" "98.32"
Set-TextValue "E51" "  +6.19%  "
